$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51 (columns B..G): Ticker, Name, Price, Market Cap, Volume, Change (24h)
$data = @(
    @("BTC", "Bitcoin", 41732, 818651115183, 8210590456, 0.4406),
    @("ETH", "Ethereum", 2471.05, 297005877842, 4353025520, 0.1045),
    @("USDT", "Tether", 0.999293, 94869792209, 12333106519, -0.01366),
    @("BNB", "BNB", 318.37, 49039706987, 338266055, 1.21854),
    @("SOL", "Solana", 93.39, 40357715765, 1048269805, 1.60475),
    @("XRP", "XRP", 0.5523400000000001, 30020044405, 575769808, 0.60766),
    @("USDC", "USDC", 0.999718, 25765376748, 2256365761, 0.01271),
    @("STETH", "Lido Staked Ether", 2471.51, 23065655266, 14200643, 0.16882),
    @("ADA", "Cardano", 0.516809, 18149997115, 256075987, 0.36765),
    @("DOGE", "Dogecoin", 0.086608, 12339343677, 1360393172, 9.10487),
    @("AVAX", "Avalanche", 33.37, 12237345227, 333436050, 2.8987),
    @("TRX", "TRON", 0.110777, 9769610452, 181906240, 0.63397),
    @("DOT", "Polkadot", 6.92, 9135382414, 125626222, 0.975),
    @("LINK", "Chainlink", 15.73, 8955172366, 435907981, -1.86797),
    @("TON", "Toncoin", 2.25, 7786165949, 33127241, 1.02866),
    @("MATIC", "Polygon", 0.7937920000000001, 7375343517, 275391293, 2.23151),
    @("WBTC", "Wrapped Bitcoin", 41670, 6594693649, 124647440, 0.36957),
    @("SHIB", "Shiba Inu", 0.00000948, 5591248068, 135244855, 0.5659999999999999),
    @("LTC", "Litecoin", 71.16, 5276264888, 237949064, 0.08298),
    @("DAI", "Dai", 0.998804, 5169190138, 73015812, -0.08746),
    @("ICP", "Internet Computer", 11.29, 5158716104, 83331325, 1.59238),
    @("UNI", "Uniswap", 6.47, 4872321380, 60082263, -0.63605),
    @("BCH", "Bitcoin Cash", 239.49, 4700144546, 112819325, 1.21861),
    @("ATOM", "Cosmos Hub", 9.789999999999999, 3758991054, 78750916, 1.25138),
    @("LEO", "LEO Token", 4.01, 3702250938, 875768, 0.91174),
    @("ETC", "Ethereum Classic", 24.69, 3539463034, 168206452, -0.38848),
    @("XLM", "Stellar", 0.116119, 3294271351, 48988246, 1.82774),
    @("OKB", "OKB", 53.79, 3226646118, 3652480, 0.89825),
    @("INJ", "Injective", 36.05, 3192947532, 93781573, 2.10215),
    @("OP", "Optimism", 3.23, 3092955273, 111087243, 2.42895),
    @("NEAR", "NEAR Protocol", 3, 3085812680, 107933666, 2.53166),
    @("XMR", "Monero", 160.21, 2901794703, 67801209, 2.85263),
    @("APT", "Aptos", 8.58, 2883779617, 82526323, 1.36963),
    @("TIA", "Celestia", 17.49, 2778428591, 76335809, 1.51008),
    @("FIL", "Filecoin", 5.51, 2739946357, 114698454, 1.00168),
    @("FDUSD", "First Digital USD", 1.001, 2616525084, 1393055429, -0.13301),
    @("LDO", "Lido DAO", 2.92, 2599157514, 54804017, 1.359),
    @("HBAR", "Hedera", 0.07657899999999999, 2579730768, 27486604, 0.84907),
    @("IMX", "Immutable", 1.95, 2574695093, 37602320, 1.8173),
    @("ARB", "Arbitrum", 1.87, 2382085324, 375978236, 4.74167),
    @("KAS", "Kaspa", 0.103367, 2317510587, 9677223, 0.3001),
    @("CRO", "Cronos", 0.083202, 2203510605, 5581298, 0.6094000000000001),
    @("MNT", "Mantle", 0.677505, 2143850275, 27161907, 0.42289),
    @("STX", "Stacks", 1.48, 2121913033, 33888053, 2.2098),
    @("VET", "VeChain", 0.02851196, 2073132010, 29553106, 0.79989),
    @("TUSD", "TrueUSD", 0.987078, 1880352089, 53543876, 0.03482),
    @("MKR", "Maker", 1999.51, 1842776132, 57513014, 3.02654),
    @("QNT", "Quant", 114.65, 1666554225, 15555296, 1.23717),
    @("SEI", "Sei", 0.678054, 1647573793, 127378418, 0.94889),
    @("TAO", "Bittensor", 264.48, 1602743893, 2831357, 8.05456)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value2 = $vals[0]
    $ws.Cells.Item($row, 3).Value2 = $vals[1]
    $ws.Cells.Item($row, 4).Value2 = $vals[2]
    $ws.Cells.Item($row, 5).Value2 = $vals[3]
    $ws.Cells.Item($row, 6).Value2 = $vals[4]
    $ws.Cells.Item($row, 7).Value2 = $vals[5]
}
